$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 209.85715
$ws.Range("I6").Value = 209.85715
$ws.Range("K6").Value = 629.5714499999999
$ws.Range("M6").Value = -517.5714499999999
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").Value = $null
$ws.Range("H28").Value = 891.5
$ws.Range("I28").Value = 910.2
$ws.Range("K28").Value = 910.2
$ws.Range("M28").Value = -425.2
$ws.Range("H43").Value = 5832
$ws.Range("I43").Value = 4998.25
$ws.Range("J43").Value = 7499.5
$ws.Range("K43").Value = 4998.25
$ws.Range("L43").Value = 7499.5
$ws.Range("M43").Value = -4929.25
$ws.Range("N43").Value = -7637.5
$ws.Range("H51").Value = 8333
$ws.Range("I51").Value = 7500
$ws.Range("J51").Value = 9999
$ws.Range("K51").Value = 7500
$ws.Range("L51").Value = 9999
$ws.Range("M51").Value = -7016
$ws.Range("N51").Value = -10967
$ws.Range("H58").Value = 350
$ws.Range("I58").Value = 350
$ws.Range("K58").Value = 1050
$ws.Range("M58").Value = -900
$ws.Range("H62").Value = 1997.5
$ws.Range("I62").Value = 1997.5
$ws.Range("K62").Value = 1997.5
$ws.Range("M62").Value = -1373.5
$ws.Range("H65").Value = 1997.5
$ws.Range("I65").Value = 1997.5
$ws.Range("K65").Value = 9987.5
$ws.Range("M65").Value = -6867.5
$ws.Range("H111").Value = 2354
$ws.Range("I111").Value = 1999
$ws.Range("K111").Value = 5997
$ws.Range("M111").Value = -2930
$ws.Range("H129").Value = 899.6667
$ws.Range("I129").Value = 899.6667
$ws.Range("K129").Value = 2699.0001
$ws.Range("M129").Value = 2300.9999
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2123
$ws.Range("I74").Value = 2060.3572
$ws.Range("J74").Value = 3000
$ws.Range("K74").Value = 2060.3572
$ws.Range("L74").Value = 3000
$ws.Range("M74").Value = -1186.3572
$ws.Range("N74").Value = -4748
$ws.Range("H77").Value = 2123
$ws.Range("I77").Value = 2060.3572
$ws.Range("J77").Value = 3000
$ws.Range("K77").Value = 10301.786
$ws.Range("L77").Value = 15000
$ws.Range("M77").Value = -5933.786
$ws.Range("N77").Value = -23736

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 830.0909
$ws.Range("I20").Value = 793.44446
$ws.Range("J20").Value = 995
$ws.Range("K20").Value = 793.44446
$ws.Range("L20").Value = 995
$ws.Range("M20").Value = -546.44446
$ws.Range("N20").Value = -1489
$ws.Range("H22").Value = 460.83334
$ws.Range("I22").Value = 460.83334
$ws.Range("K22").Value = 460.83334
$ws.Range("M22").Value = -287.83334
$ws.Range("H134").Value = 1205.2307
$ws.Range("I134").Value = 1205.2307
$ws.Range("K134").Value = 3615.6921
$ws.Range("M134").Value = -1080.6921

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 5948.364
$ws.Range("I22").Value = 7103.5557
$ws.Range("K22").Value = 7103.5557
$ws.Range("M22").Value = -6753.5557
$ws.Range("H26").Value = 600
$ws.Range("I26").Value = 600
$ws.Range("K26").Value = 600
$ws.Range("M26").Value = -313
$ws.Range("H31").Value = 1527.4
$ws.Range("I31").Value = 1499.3334
$ws.Range("J31").Value = 1569.5
$ws.Range("K31").Value = 1499.3334
$ws.Range("L31").Value = 1569.5
$ws.Range("M31").Value = -1204.3334
$ws.Range("N31").Value = -2159.5
$ws.Range("H34").Value = 1527.4
$ws.Range("I34").Value = 1499.3334
$ws.Range("J34").Value = 1569.5
$ws.Range("K34").Value = 1499.3334
$ws.Range("L34").Value = 1569.5
$ws.Range("M34").Value = -1297.3334
$ws.Range("N34").Value = -1973.5
$ws.Range("H42").Value = 30000
$ws.Range("I42").Value = 30000
$ws.Range("K42").Value = 30000
$ws.Range("M42").Value = -29407
$ws.Range("H55").Value = 40000
$ws.Range("I55").Value = 40000
$ws.Range("K55").Value = 40000
$ws.Range("M55").Value = -39685
$ws.Range("H58").Value = 2958.1
$ws.Range("I58").Value = 3174
$ws.Range("J58").Value = 2094.5
$ws.Range("K58").Value = 3174
$ws.Range("L58").Value = 2094.5
$ws.Range("M58").Value = -2971
$ws.Range("N58").Value = -2500.5
$ws.Range("H86").Value = 600602.6
$ws.Range("I86").Value = 749503.25
$ws.Range("K86").Value = 749503.25
$ws.Range("M86").Value = -748380.25
$ws.Range("H89").Value = 600602.6
$ws.Range("I89").Value = 749503.25
$ws.Range("K89").Value = 3747516.25
$ws.Range("M89").Value = -3741900.25
$ws.Range("H136").Value = 2958.1
$ws.Range("I136").Value = 3174
$ws.Range("J136").Value = 2094.5
$ws.Range("K136").Value = 9522
$ws.Range("L136").Value = 6283.5
$ws.Range("M136").Value = -6972
$ws.Range("N136").Value = -11383.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 300
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 300
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 900
$ws.Range("M7").Value = $null
$ws.Range("N7").Value = -1124
$ws.Range("H68").Value = 5000
$ws.Range("J68").Value = 5000
$ws.Range("L68").Value = 15000
$ws.Range("N68").Value = -16622
$ws.Range("H71").Value = 5000
$ws.Range("J71").Value = 5000
$ws.Range("L71").Value = 45000
$ws.Range("N71").Value = -53112
$ws.Range("H107").Value = 1035.4286
$ws.Range("J107").Value = 1035.4286
$ws.Range("L107").Value = 3106.2858
$ws.Range("N107").Value = -6946.2858
$ws.Range("H112").Value = 2181.6667
$ws.Range("I112").Value = 1955.25
$ws.Range("J112").Value = 2634.5
$ws.Range("K112").Value = 5865.75
$ws.Range("L112").Value = 7903.5
$ws.Range("M112").Value = -4757.75
$ws.Range("N112").Value = -10119.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 36666.5
$ws.Range("I5").Value = 10000
$ws.Range("J5").Value = 49999.75
$ws.Range("K5").Value = 10000
$ws.Range("L5").Value = 49999.75
$ws.Range("M5").Value = -9888
$ws.Range("N5").Value = -50223.75
$ws.Range("H122").Value = 2799
$ws.Range("I122").Value = 2799
$ws.Range("K122").Value = 8397
$ws.Range("M122").Value = -5947
$ws.Range("H132").Value = 3999.6667
$ws.Range("I132").Value = 3999.6667
$ws.Range("K132").Value = 11999.0001
$ws.Range("M132").Value = -9469.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4247.5
$ws.Range("J40").Value = 5000
$ws.Range("L40").Value = 5000
$ws.Range("N40").Value = -5272
$ws.Range("H61").Value = 705.0909
$ws.Range("I61").Value = 705.0909
$ws.Range("K61").Value = 705.0909
$ws.Range("M61").Value = -503.0909
$ws.Range("H113").Value = 705.0909
$ws.Range("I113").Value = 705.0909
$ws.Range("K113").Value = 705.0909
$ws.Range("M113").Value = 1464.9091

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 2999.5
$ws.Range("I3").Value = 3999
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 3999
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -3885
$ws.Range("N3").Value = -2228
$ws.Range("H113").Value = 443.4
$ws.Range("I113").Value = 427.45456
$ws.Range("J113").Value = 462.8889
$ws.Range("K113").Value = 1282.36368
$ws.Range("L113").Value = 1388.6667
$ws.Range("M113").Value = 887.6363200000001
$ws.Range("N113").Value = -5728.6667
$ws.Range("H141").Value = 90000
$ws.Range("J141").Value = 90000
$ws.Range("L141").Value = 90000
$ws.Range("N141").Value = -100360

Write-Host "Applied all cell updates"
